$d = $word.ActiveDocument

# The document has the Pearson logo (descr contains "PearsonLogo.png")
# embedded in every footer as "image1.png" and the BTec logo (descr
# "BTec_Logo-Orange") embedded in every header as "image2.jpg". This
# edit swaps those two picture names: the Pearson logo becomes
# "image2.png" and the BTec logo becomes "image1.jpg" (the underlying
# media part / relationship target is untouched - only the shape's
# display name changes).
#
# InlineShape has no settable Name property (matches real Word's object
# model - only Shape/ShapeRange expose Name), so each picture is
# temporarily promoted to a floating Shape, renamed, then converted back
# to an inline picture so the <wp:inline> layout is preserved.
function Rename-Logo($inlineShape, $newName) {
    $shape = $inlineShape.ConvertToShape()
    $shape.Name = $newName
    [void]$shape.ConvertToInlineShape()
}

for ($s = 1; $s -le $d.Sections.Count; $s++) {
    $sec = $d.Sections.Item($s)

    for ($i = 1; $i -le $sec.Headers.Count; $i++) {
        $hdr = $sec.Headers.Item($i)
        if ($hdr.Exists) {
            for ($j = 1; $j -le $hdr.Range.InlineShapes.Count; $j++) {
                $ishp = $hdr.Range.InlineShapes.Item($j)
                if ($ishp.AlternativeText -eq "BTec_Logo-Orange") {
                    Rename-Logo $ishp "image1.jpg"
                }
            }
        }
    }

    for ($i = 1; $i -le $sec.Footers.Count; $i++) {
        $ftr = $sec.Footers.Item($i)
        if ($ftr.Exists) {
            for ($j = 1; $j -le $ftr.Range.InlineShapes.Count; $j++) {
                $ishp = $ftr.Range.InlineShapes.Item($j)
                if ($ishp.AlternativeText -like "*PearsonLogo.png") {
                    Rename-Logo $ishp "image2.png"
                }
            }
        }
    }
}

Write-Output "done"
